# Add a new "capacitors" worksheet after the existing "Sheet1", populate it
# with the chosen capacitor's data, and make it the active/selected sheet
# (mirrors what happened when the author picked capacitor types in Excel).

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws.Name = "capacitors"

# Header row.
$ws.Range("A1").Value = "Value"
$ws.Range("B1").Value = "size"
$ws.Range("C1").Value = "Voltage"
$ws.Range("D1").Value = "link"

# Data row - write the link before the other values so new shared strings
# are interned in the same order as the source workbook.
$ws.Range("D2").Value = "https://www.hestore.hu/prod_10024356.html"
$ws.Range("A2").Value = "100nf"
$ws.Range("B2").Value = 603
$ws.Range("C2").Value = "50V"

# Make "capacitors" the active sheet/tab, with its own zoom + selection.
$ws.Activate()
$excel.ActiveWindow.Zoom = 220
[void]$ws.Range("C6").Select()
